$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.178.94"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "3.554.60"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'603.92"
$ws.Range("E5").Value = "  -0.50%  "

# Row 6
$ws.Range("D6").Value = "'143.29"
$ws.Range("E6").Value = "  -1.13%  "

# Row 7
$ws.Range("D7").Value = "3.553.05"
$ws.Range("E7").Value = "  +1.16%  "

# Row 9
$ws.Range("D9").Value = "'0.489"

# Row 10
$ws.Range("E10").Value = "  -0.62%  "

# Row 11
$ws.Range("E11").Value = "  -1.86%  "

# Row 12
$ws.Range("E12").Value = "  -0.43%  "

# Row 13
$ws.Range("D13").Value = "4.158.96"
$ws.Range("E13").Value = "  +1.26%  "

# Row 14
$ws.Range("E14").Value = "  -1.30%  "

# Row 15
$ws.Range("D15").Value = "'29.99"
$ws.Range("E15").Value = "  -1.18%  "

# Row 16
$ws.Range("D16").Value = "3.568.73"
$ws.Range("E16").Value = "  +1.72%  "

# Row 17
$ws.Range("D17").Value = "66.239.68"
$ws.Range("E17").Value = "  -0.13%  "

# Row 18
$ws.Range("E18").Value = "  -0.60%  "

# Row 19
$ws.Range("D19").Value = "'11.37"
$ws.Range("E19").Value = "  +8.84%  "

# Row 20
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("E21").Value = "  -1.49%  "

# Row 22
$ws.Range("D22").Value = "'428.61"
$ws.Range("E22").Value = "  +0.49%  "

# Row 23
$ws.Range("D23").Value = "'0.607"
$ws.Range("E23").Value = "  +1.70%  "

# Row 24
$ws.Range("D24").Value = "'79.71"
$ws.Range("E24").Value = "  +2.19%  "

# Row 25
$ws.Range("D25").Value = "3.698.16"
$ws.Range("E25").Value = "  +1.54%  "

# Row 26
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("E27").Value = "  -3.78%  "

# Row 28
$ws.Range("E28").Value = "  +0.49%  "

# Row 29
$ws.Range("D29").Value = "'9.06"
$ws.Range("E29").Value = "  -2.64%  "

# Row 30
$ws.Range("D30").Value = "'7.83"
$ws.Range("E30").Value = "  -2.01%  "

# Row 31
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("D32").Value = "3.552.69"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("E33").Value = "  +0.48%  "

# Row 34
$ws.Range("E34").Value = "  -2.23%  "

# Row 35
$ws.Range("E35").Value = "  -9.32%  "

# Row 37
$ws.Range("D37").Value = "'7.79"
$ws.Range("E37").Value = "  -0.11%  "

# Row 38
$ws.Range("E38").Value = "  -1.41%  "

# Row 39
$ws.Range("D39").Value = "'5.53"
$ws.Range("E39").Value = "  -1.54%  "

# Row 40
$ws.Range("D40").Value = "'173.97"
$ws.Range("E40").Value = "  +1.95%  "

# Row 41
$ws.Range("D41").Value = "'0.0846"
$ws.Range("E41").Value = "  -1.65%  "

# Row 42
$ws.Range("D42").Value = "'5.18"
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("E43").Value = "  +0.14%  "

# Row 44
$ws.Range("D44").Value = "'1.92"
$ws.Range("E44").Value = "  +1.03%  "

# Row 45
$ws.Range("E45").Value = "  +1.21%  "

# Row 46
$ws.Range("E46").Value = "  -0.02%  "

# Row 47
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("D48").Value = "'24.79"
$ws.Range("E48").Value = "  -4.19%  "

# Row 49
$ws.Range("E49").Value = "  -2.08%  "

# Row 50
$ws.Range("E50").Value = "  -1.11%  "

# Row 51
$ws.Range("D51").Value = "'22.95"
$ws.Range("E51").Value = "  +1.27%  "
